$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$shp.Table.ApplyStyle("{90BAF084-68BA-4B18-A833-53D3B7BAAED8}")
